# Commit: "Remove form_id from remaining forms"
#
# The "settings" sheet has columns: form_title | form_id | version | style | namespaces
# form_id (and its value "case_event") is no longer needed, so the whole
# column is removed and everything to its right shifts one column left:
#   form_title | version | style | namespaces

$wb       = $excel.ActiveWorkbook
$survey   = $wb.Worksheets.Item("survey")
$settings = $wb.Worksheets.Item("settings")

# ---------------------------------------------------------------------------
# 1. Cell comments are anchored to fixed cells and do NOT follow a column
#    shift automatically, so grab their current text now (while B1..E1 still
#    point at the original form_id/version/style/namespaces comments) and
#    re-apply it to the right cells after the column is removed:
#      old B1 (form_id comment)     -> removed entirely
#      old C1 (version comment)     -> becomes new B1
#      old D1 (style/pages comment) -> becomes new C1
#      old E1 (namespaces comment)  -> becomes new D1
# ---------------------------------------------------------------------------
$versionText    = $settings.Range("C1").Comment.Text()
$pagesText      = $settings.Range("D1").Comment.Text()
$namespacesText = $settings.Range("E1").Comment.Text()

# ---------------------------------------------------------------------------
# 2. Delete the form_id column (column B). Cell values, shared strings,
#    formulas and column widths to the right all shift left automatically.
# ---------------------------------------------------------------------------
$settings.Columns.Item(2).Delete()

# ---------------------------------------------------------------------------
# 3. Re-point the surviving comments and drop the now-empty E1 comment.
# ---------------------------------------------------------------------------
[void]$settings.Range("B1").Comment.Text($versionText)
[void]$settings.Range("C1").Comment.Text($pagesText)
[void]$settings.Range("D1").Comment.Text($namespacesText)
$settings.Range("E1").Comment.Delete()

# ---------------------------------------------------------------------------
# 4. survey sheet: the conditional formatting had accumulated a leftover,
#    row-27-specific split (an artifact of an earlier row insert/delete).
#    Collapse it back down to clean, contiguous ranges.
# ---------------------------------------------------------------------------
$cf = $survey.Range("A1:Z30").FormatConditions
$cf.Item(1).ModifyAppliesToRange($survey.Range("A2:D10000"))
$cf.Item(7).ModifyAppliesToRange($survey.Range("C2:C10000"))
for ($i = 13; $i -ge 8; $i--) {
    $cf.Item($i).Delete()
}

# ---------------------------------------------------------------------------
# 5. Restore the selections / active sheet as left by the author.
# ---------------------------------------------------------------------------
[void]$settings.Range("B1").Select()
[void]$survey.Activate()
[void]$survey.Range("A6").Select()
